$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 123: new expense entry (19 Nov 2024, dinner) ---
$ws.Range("A122:F122").Copy()
$ws.Range("A123:F123").PasteSpecial(-4122)
$ws.Range("B123").Value = 45615
$ws.Range("C123").Value = "晚饭"
$ws.Range("D123").Value = -18
$ws.Range("E123").Value = "麻辣烫（13）+烤肠（5）"
$ws.Range("F123").Formula = "=F122+D123"

# --- Row 124: daily summary row (小结) for 19 Nov 2024 ---
$ws.Range("A121:H121").Copy()
$ws.Range("A124:H124").PasteSpecial(-4122)
$ws.Range("B124").Value = 45615
$ws.Range("C124").Value = "小结"
$ws.Range("D124").Formula = "=SUM(D122:D123)"
$ws.Range("E124").Value = "*"
$ws.Range("F124").Value = 413.45000000000027

# --- Row 125: new expense entry (20 Nov 2024, lunch) ---
$ws.Range("A122:F122").Copy()
$ws.Range("A125:F125").PasteSpecial(-4122)
$ws.Range("B125").Value = 45616
$ws.Range("C125").Value = "午饭"
$ws.Range("D125").Value = -3.2
$ws.Range("E125").Value = "两份六两大米"
$ws.Range("F125").Formula = "=F124+D125"

$ws.Range("F121").Select()
